$d = $word.ActiveDocument

# --- Edit 1: "At Least Once System" -> "At Least Once System," -----------
# (the bold run "At Least Once System" becomes "At Least Once System,")
$r1 = $d.Content
$found1 = $r1.Find.Execute("At Least Once System", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $rInsert1 = $d.Range($r1.End, $r1.End)
    $rInsert1.InsertBefore(",")
}

# --- Edit 2: bold "enable.idempotent=true" inside the sentence -----------
# " setting enable.idempotent=true will help you ensure..." ->
# " setting " + bold("enable.idempotent=true") + " will help you ensure..."
$r2 = $d.Content
$found2 = $r2.Find.Execute("enable.idempotent=true", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $rBold = $d.Range($r2.Start, $r2.End)
    $rBold.Font.Bold = $true
}
